$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New coin rankings/prices/volumes refreshed from the source feed.
# Column D holds numeric-looking price strings (e.g. "1.000", "27.212.74")
# that must stay literal text, so it is pre-formatted as Text before the
# values are written (only for the rows whose price actually changes).
$changedPriceRows = @(
2, 3, 5, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39, 40, 41, 42, 43, 44, 45, 46, 47, 48, 49, 50, 51
)
foreach ($r in $changedPriceRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$data = @(
    @($null, $null, "27.212.74", "  +0.62%  "),
    @($null, $null, "1.903.03", "  +0.55%  "),
    @($null, $null, $null, "  -0.20%  "),
    @($null, $null, "305.92", "  -0.23%  "),
    @($null, $null, $null, "  -0.22%  "),
    @($null, $null, "0.5403", "  +3.48%  "),
    @($null, $null, "0.3803", "  +1.20%  "),
    @($null, $null, "0.07288", $null),
    @($null, $null, "22.04", "  +4.59%  "),
    @($null, $null, "0.9020", "  +0.44%  "),
    @($null, $null, "0.08182", "  +0.26%  "),
    @($null, $null, "95.54", "  -0.66%  "),
    @($null, $null, "5.345", "  +0.87%  "),
    @($null, $null, "0.9980", "  -0.48%  "),
    @($null, $null, "14.80", "  +1.63%  "),
    @($null, $null, "0.000008632", "  +0.59%  "),
    @($null, $null, "1.000", "  -0.27%  "),
    @("WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "27.254.17", "  +0.67%  "),
    @("Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "5.044", "  -0.65%  "),
    @("Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "10.81", "  +1.18%  "),
    @("Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "6.509", "  +1.59%  "),
    @("Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "148.22", "  -0.34%  "),
    @("LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "2.304", "  +0.67%  "),
    @("EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "18.34", "  +0.95%  "),
    @("Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "1.755", "  +1.18%  "),
    @("BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "116.62", "  +1.40%  "),
    @("InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "4.853", "  +1.50%  "),
    @("Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "4.650", "  -4.07%  "),
    @("Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.09204", "  -0.17%  "),
    @("ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "0.8238", "  +4.57%  "),
    @("Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.05059", "  +0.60%  "),
    @("ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "1.220", "  +0.86%  "),
    @("HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "3.011", "  +1.15%  "),
    @("MXToken", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx", "3.316", "  -3.10%  "),
    @("RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "2.699", "  +3.76%  "),
    @("TheSandbox", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand", "0.6007", "  +5.32%  "),
    @("VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.02000", "  +0.72%  "),
    @("TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "1.074", "  -0.02%  "),
    @("Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "9.249", "  +2.38%  "),
    @("FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "6.649", "  +1.60%  "),
    @("Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "115.91", "  -0.16%  "),
    @("Decentraland", "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana", "0.5155", "  +6.18%  "),
    @("Algorand", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo", "0.1529", "  +1.01%  "),
    @("PaxDollar", "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp", "0.9994", "  -0.30%  "),
    @($null, $null, "10.13", "  +0.94%  "),
    @("NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "1.637", "  +0.93%  "),
    @("Elrond", "https://coinranking.com/coin/omwkOTglq+elrond-egld", "38.07", "  -0.25%  "),
    @("Cronos", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro", "0.06097", "  +2.85%  "),
    @("Aave", "https://coinranking.com/coin/ixgUfzmLR+aave-aave", "63.44", "  -0.02%  "),
    @("EOS", "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos", "0.9223", "  +1.06%  ")
)

$row = 2
foreach ($item in $data) {
    if ($null -ne $item[0]) { $ws.Cells.Item($row, 2).Value = $item[0] }
    if ($null -ne $item[1]) { $ws.Cells.Item($row, 3).Value = $item[1] }
    if ($null -ne $item[2]) { $ws.Cells.Item($row, 4).Value = $item[2] }
    if ($null -ne $item[3]) { $ws.Cells.Item($row, 5).Value = $item[3] }
    $row++
}
